$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (was "UniformA-HW50.xpc" -> "UniformA")
$ws.Name = "UniformA"

# Append a new data row (row 16), mirroring the pattern of the existing
# rows (A = index, B = shared label text reused from row 15, C:P = 1s).
# Copy formatting from A15 (bold/border/center style) onto A16 first.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16:P16").Value = 1
